$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'reflective jogging wear',
    'reflective leg band',
    'reflective leg bands',
    'reflective leggings',
    'reflective leggings women',
    'reflective overalls',
    'reflective pants',
    'reflective pants band',
    'reflective pants for work',
    'reflective pants women',
    'reflective post panel',
    'reflective rainbow leggings women',
    'reflective run gear',
    'reflective running apparel',
    'reflective running clothes',
    'reflective running clothes for women',
    'reflective running gear black',
    'reflective running gear for shoes',
    'reflective running gear for women',
    'reflective running gear women',
    'reflective running jacket women',
    'reflective running leg band',
    'reflective running leggings',
    'reflective running leggings for women',
    'reflective running lights',
    'reflective running pants',
    'reflective running pants men',
    'reflective running pants women',
    'reflective running tank',
    'reflective running tights',
    'reflective running tights womens',
    'reflective running wear',
    'reflective running women',
    'reflective string',
    'reflective tights women',
    'reflective wear',
    'reflective wear for women',
    'reflective winter gear',
    'reflective work pants',
    'reflective workout leggings',
    'reflex capri',
    'regular panty hose',
    'reinforced leggings',
    'resistance for running',
    'resistance running',
    'resistance workout',
    'riding breaches',
    'riding leggings',
    'riding leggings women',
    'riding pants for women',
    'riding pants women',
    'riding standing up',
    'riding tights',
    'riding tights for women',
    'right hip pain',
    'right side hip pain',
    'ring bombs size 7',
    'ring leggings',
    'ring power',
    'ring size chart',
    'ring sizing chart',
    'rip away pants basketball',
    'rip away pants mens',
    'rip away track pants',
    'rip muscle x',
    'ripped leggings for women',
    'ripped tights',
    'ripped tights for women',
    'rise against end game',
    'rise basketball shoes',
    'rise from hate',
    'rise jacket',
    'rise up women',
    'rise water',
    'rise woman',
    'road bike art',
    'road crew vest',
    'road leggings',
    'road runner',
    'road runner leggings',
    'road runner sports',
    'road runner sports shorts women',
    'road runners',
    'road skis',
    'roadrunner bike',
    'roadrunner game',
    'roast everything',
    'rock climbing pants for men',
    'rock climbing pants women',
    'rock climbing pants womens',
    'rock climbing thailand',
    'rock climbing xl holds',
    'rock tape kinesiology tape',
    'rock tape kinesiology tape knee',
    'rocks leggings',
    'roots clothing women',
    'roots pants',
    'ropa de compresion mujer',
    'ropa para correr',
    'ropa running mujer'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

